# Applies:
#  1. slide5: merge "per database" / ", " / "not counted toward buffer pool limit"
#     runs into a single run.
#  2. Inserts two new slides ("SQL Server Features" index slide and a blank
#     "Row-Level Security" slide) right before the final "Resources" slide.

$p = $ppt.ActivePresentation

# --- 1. Fix the run split on slide 5 ------------------------------------
$slide5 = $p.Slides.Item(5)
$bodyTr = $slide5.Shapes.Item(2).TextFrame.TextRange
$full = $bodyTr.Text
$idx = $full.IndexOf("per database")
$target = $bodyTr.Characters($idx + 1, 50)
$target.Text = "per database, not counted toward buffer pool limit"

# --- 2. Insert the new "SQL Server Features" slide (index 6) -----------
$slide6 = $p.Slides.Add(6, 2)
$slide6.Shapes.Item(1).TextFrame.TextRange.Text = "SQL Server Features"

$body6 = $slide6.Shapes.Item(2).TextFrame.TextRange
$body6.Text = "Security`rRow-Level Security`rAlways Encrypted`rUtility`rSnapshots`rFilestream`rPerformance`rColumnstore Indexes`rPartitioning`rCompression`rIn-Memory OLTP (Hekaton)"

$body6.Paragraphs(2).IndentLevel = 2
$body6.Paragraphs(3).IndentLevel = 2
$body6.Paragraphs(5).IndentLevel = 2
$body6.Paragraphs(6).IndentLevel = 2
$body6.Paragraphs(8).IndentLevel = 2
$body6.Paragraphs(9).IndentLevel = 2
$body6.Paragraphs(10).IndentLevel = 2
$body6.Paragraphs(11).IndentLevel = 2

$body6.Paragraphs(2).Font.Color.RGB = 255

# --- 3. Insert the new (still empty) "Row-Level Security" slide (index 7) -
$slide7 = $p.Slides.Add(7, 2)
$slide7.Shapes.Item(1).TextFrame.TextRange.Text = "Row-Level Security"
